$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.417.39"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "3.440.29"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'579.44"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'148.51"
$ws.Range("E6").Value = "  +8.78%  "
$ws.Range("D7").Value = "3.441.63"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'7.79"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "'0.391"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "4.028.48"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'27.98"
$ws.Range("E14").Value = "  +6.16%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "3.444.03"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "61.532.98"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "  +8.25%  "
$ws.Range("D20").Value = "'14.37"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "'9.44"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'388.53"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").Value = "'0.570"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").Value = "3.590.94"
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'72.66"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E29").Value = "  +6.84%  "
$ws.Range("D30").Value = "'7.80"
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -13.93%  "
$ws.Range("D33").Value = "'8.25"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "'2.17"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'23.97"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "'5.28"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'7.07"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "'166.01"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "'0.0792"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("D42").Value = "'26.25"
$ws.Range("E42").Value = "  +8.75%  "
$ws.Range("D43").Value = "'0.794"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "'1.72"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").Value = "2.613.87"
$ws.Range("E48").Value = "  +8.33%  "
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("D50").Value = "'7.04"
$ws.Range("E50").Value = "  +3.41%  "
$ws.Range("E51").Value = "  -1.02%  "
